# Auto-generated Excel COM-interop script
# Applies cell text updates per the target diff, forcing text storage
# (avoids Excel auto-converting numeric-looking strings into numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextCell $ws "D2" "29.312.28"
Set-TextCell $ws "E2" "  +1.11%  "
Set-TextCell $ws "D3" "1.910.14"
Set-TextCell $ws "E3" "  +1.36%  "
Set-TextCell $ws "E4" "  -0.10%  "
Set-TextCell $ws "D5" "322.18"
Set-TextCell $ws "E5" "  -2.69%  "
Set-TextCell $ws "E6" "  -0.09%  "
Set-TextCell $ws "D7" "0.4710"
Set-TextCell $ws "E7" "  +2.38%  "
Set-TextCell $ws "D8" "0.4047"
Set-TextCell $ws "E8" "  -0.35%  "
Set-TextCell $ws "B9" "Dogecoin"
Set-TextCell $ws "C9" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell $ws "D9" "0.08026"
Set-TextCell $ws "E9" "  +0.58%  "
Set-TextCell $ws "B10" "Polygon"
Set-TextCell $ws "C10" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell $ws "D10" "0.9987"
Set-TextCell $ws "E10" "  +0.88%  "
Set-TextCell $ws "B11" "Solana"
Set-TextCell $ws "C11" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell $ws "D11" "22.58"
Set-TextCell $ws "E11" "  +4.28%  "
Set-TextCell $ws "B12" "WrappedEther"
Set-TextCell $ws "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws "D12" "1.916.41"
Set-TextCell $ws "E12" "  +1.61%  "
Set-TextCell $ws "B13" "Polkadot"
Set-TextCell $ws "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws "D13" "5.869"
Set-TextCell $ws "E13" "  -0.54%  "
Set-TextCell $ws "B14" "Chainlink"
Set-TextCell $ws "C14" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws "D14" "7.091"
Set-TextCell $ws "E14" "  +0.42%  "
Set-TextCell $ws "B15" "Litecoin"
Set-TextCell $ws "C15" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell $ws "D15" "89.44"
Set-TextCell $ws "E15" "  +1.14%  "
Set-TextCell $ws "B16" "BinanceUSD"
Set-TextCell $ws "C16" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell $ws "D16" "1.001"
Set-TextCell $ws "E16" "  -0.14%  "
Set-TextCell $ws "B17" "TRON"
Set-TextCell $ws "C17" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell $ws "D17" "0.06617"
Set-TextCell $ws "E17" "  +0.86%  "
Set-TextCell $ws "B18" "ShibaInu"
Set-TextCell $ws "C18" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell $ws "D18" "0.00001028"
Set-TextCell $ws "E18" "  -0.30%  "
Set-TextCell $ws "B19" "Avalanche"
Set-TextCell $ws "C19" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell $ws "D19" "17.60"
Set-TextCell $ws "E19" "  +1.02%  "
Set-TextCell $ws "B20" "Dai"
Set-TextCell $ws "C20" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws "D20" "1.001"
Set-TextCell $ws "E20" "  -0.06%  "
Set-TextCell $ws "B21" "WrappedBTC"
Set-TextCell $ws "C21" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell $ws "D21" "29.326.27"
Set-TextCell $ws "E21" "  +1.02%  "
Set-TextCell $ws "B22" "Uniswap"
Set-TextCell $ws "C22" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell $ws "D22" "5.515"
Set-TextCell $ws "E22" "  +1.78%  "
Set-TextCell $ws "B23" "Cosmos"
Set-TextCell $ws "C23" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell $ws "D23" "11.41"
Set-TextCell $ws "E23" "  -0.01%  "
Set-TextCell $ws "B24" "Toncoin"
Set-TextCell $ws "C24" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws "D24" "2.201"
Set-TextCell $ws "E24" "  -0.33%  "
Set-TextCell $ws "B25" "WrappedliquidstakedEther2.0"
Set-TextCell $ws "C25" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell $ws "D25" "2.163.60"
Set-TextCell $ws "E25" "  +2.51%  "
Set-TextCell $ws "B26" "Monero"
Set-TextCell $ws "C26" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws "D26" "154.40"
Set-TextCell $ws "E26" "  -1.58%  "
Set-TextCell $ws "B27" "EthereumClassic"
Set-TextCell $ws "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws "D27" "19.73"
Set-TextCell $ws "E27" "  +0.76%  "
Set-TextCell $ws "B28" "InternetComputer(DFINITY)"
Set-TextCell $ws "C28" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws "D28" "6.016"
Set-TextCell $ws "E28" "  +9.60%  "
Set-TextCell $ws "B29" "LidoDAOToken"
Set-TextCell $ws "C29" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell $ws "D29" "2.099"
Set-TextCell $ws "E29" "  +0.44%  "
Set-TextCell $ws "B30" "BitcoinCash"
Set-TextCell $ws "C30" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell $ws "D30" "117.74"
Set-TextCell $ws "E30" "  +0.23%  "
Set-TextCell $ws "B31" "ImmutableX"
Set-TextCell $ws "C31" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws "D31" "1.070"
Set-TextCell $ws "E31" "  +5.51%  "
Set-TextCell $ws "B32" "Stellar"
Set-TextCell $ws "C32" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws "D32" "0.09487"
Set-TextCell $ws "E32" "  +1.71%  "
Set-TextCell $ws "B33" "ARBITRUM"
Set-TextCell $ws "C33" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws "D33" "1.414"
Set-TextCell $ws "E33" "  +0.79%  "
Set-TextCell $ws "B34" "HuobiToken"
Set-TextCell $ws "C34" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell $ws "D34" "3.541"
Set-TextCell $ws "E34" "  -1.70%  "
Set-TextCell $ws "B35" "Filecoin"
Set-TextCell $ws "C35" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws "D35" "5.360"
Set-TextCell $ws "E35" "  +1.63%  "
Set-TextCell $ws "B36" "Hedera"
Set-TextCell $ws "C36" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws "D36" "0.06059"
Set-TextCell $ws "E36" "  +0.01%  "
Set-TextCell $ws "B37" "VeChain"
Set-TextCell $ws "C37" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws "D37" "0.02245"
Set-TextCell $ws "E37" "  +1.15%  "
Set-TextCell $ws "B38" "FraxShare"
Set-TextCell $ws "C38" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws "D38" "8.186"
Set-TextCell $ws "E38" "  -1.18%  "
Set-TextCell $ws "B39" "TrustWalletToken"
Set-TextCell $ws "C39" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws "D39" "1.177"
Set-TextCell $ws "E39" "  +0.60%  "
Set-TextCell $ws "B40" "TheSandbox"
Set-TextCell $ws "C40" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell $ws "D40" "0.5831"
Set-TextCell $ws "E40" "  +0.85%  "
Set-TextCell $ws "B41" "RenderToken"
Set-TextCell $ws "C41" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws "D41" "2.493"
Set-TextCell $ws "E41" "  +10.40%  "
Set-TextCell $ws "D42" "0.1833"
Set-TextCell $ws "E42" "  +0.46%  "
Set-TextCell $ws "B43" "Aptos"
Set-TextCell $ws "C43" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell $ws "D43" "10.09"
Set-TextCell $ws "E43" "  -0.26%  "
Set-TextCell $ws "B44" "Cronos"
Set-TextCell $ws "C44" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws "D44" "0.07875"
Set-TextCell $ws "E44" "  +6.09%  "
Set-TextCell $ws "B45" "WEMIXToken"
Set-TextCell $ws "C45" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws "D45" "1.277"
Set-TextCell $ws "E45" "  +1.11%  "
Set-TextCell $ws "B46" "Decentraland"
Set-TextCell $ws "C46" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell $ws "D46" "0.5494"
Set-TextCell $ws "E46" "  +0.82%  "
Set-TextCell $ws "B47" "EnergySwap"
Set-TextCell $ws "C47" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws "D47" "12.04"
Set-TextCell $ws "E47" "  +0.45%  "
Set-TextCell $ws "B48" "NEARProtocol"
Set-TextCell $ws "C48" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws "D48" "1.914"
Set-TextCell $ws "E48" "  +0.89%  "
Set-TextCell $ws "B49" "Quant"
Set-TextCell $ws "C49" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell $ws "D49" "113.16"
Set-TextCell $ws "E49" "  +1.97%  "
Set-TextCell $ws "B50" "Elrond"
Set-TextCell $ws "C50" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextCell $ws "D50" "44.20"
Set-TextCell $ws "E50" "  -3.32%  "
Set-TextCell $ws "B51" "WOONetwork"
Set-TextCell $ws "C51" "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextCell $ws "D51" "0.2909"
Set-TextCell $ws "E51" "  +3.61%  "
